$d = $word.ActiveDocument

function Find-ParaByPrefix($prefix) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $cand = $d.Paragraphs($i)
        if ($cand.Range.Text.StartsWith($prefix)) {
            return $cand
        }
    }
    throw ("paragraph starting with '" + $prefix + "' not found")
}

function Set-ParaRuns($para, $innerXml) {
    $full = $para.Range
    # Exclude the trailing paragraph mark so the <w:p> wrapper (and its
    # paraId/textId/rsid attributes plus <w:pPr>) is left untouched; only
    # the run content inside the paragraph gets replaced.
    $r = $d.Range($full.Start, $full.End - 1)
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $innerXml + '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg)
}

# "Navigate to QuoteAggregatorApi and run the following command: docker-compose up"
# -> mark "QuoteAggregatorApi" with spell-check proofing runs.
$paraApi = Find-ParaByPrefix "Navigate to QuoteAggregatorApi"
$xmlApi = '<w:r><w:t xml:space="preserve">Navigate to </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>QuoteAggregatorApi</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> and run the following command: docker-compose up</w:t></w:r>'
Set-ParaRuns $paraApi $xmlApi

# "If you have maven installed, you can run mvn spring-boot:run in the same folder"
# -> mark "mvn" and "spring-boot:run" with spell-check proofing runs.
$paraMaven = Find-ParaByPrefix "If you have maven installed"
$xmlMaven = '<w:r><w:t xml:space="preserve">If you have maven installed, you can run </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>mvn</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>spring-boot:run</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> in the same folder</w:t></w:r>'
Set-ParaRuns $paraMaven $xmlMaven

# "Navigate to QuoteAggregatorApp, and run the following command: npm run dev"
# -> fix the path to "QuoteAggregatorApp/QuoteAggregatorFrontEnd" and mark the
#    non-dictionary words with spell-check proofing runs.
$paraApp = Find-ParaByPrefix "Navigate to QuoteAggregatorApp"
$xmlApp = '<w:r><w:t xml:space="preserve">Navigate to </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>QuoteAggregatorApp</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>/</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>QuoteAggregator</w:t></w:r>' +
    '<w:r><w:t>FrontEnd</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">, and run the following command: </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>npm</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> run dev</w:t></w:r>'
Set-ParaRuns $paraApp $xmlApp
